$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title in A1 with the new migration wave date
$ws.Range("A1").Value = "Hotcarding Spreadsheet - Migration Wave 2025-10-16 00:00:00"

# Append the new migration wave row (row 6)
$ws.Range("A6").Value = "2025-10-16 00:00:00"
$ws.Range("B6").Value = "YYY"
$ws.Range("C6").Value = "123ABX007"
$ws.Range("D6").Value = "FISB"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "PaymentsOne Debit"
$ws.Range("G6").Value = "Basic"
$ws.Range("H6").Value = "Offshore"
$ws.Range("I6").Value = "NA"
$ws.Range("J6").Value = "NA"

# The source rows (4, 5) carry no explicit cell style even though their
# columns define a default style; reset the new row the same way so it
# doesn't pick up the column's style index.
$ws.Range("A6:J6").Style = "Normal"
